$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.153.32"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "'2.421.76"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'563.12"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("D6").Value = "'144.21"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'0.532"
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("D9").Value = "'2.418.77"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").Value = "'26.11"
$ws.Range("E14").Value = "  +3.34%  "
$ws.Range("E15").Value = "  +5.43%  "
$ws.Range("D16").Value = "'2.860.00"
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").Value = "'61.998.58"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "'2.421.55"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").Value = "'11.19"
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("D20").Value = "'324.59"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").Value = "'4.18"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'65.40"
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").Value = "'9.02"
$ws.Range("E26").Value = "  +4.67%  "
$ws.Range("D27").Value = "'586.76"
$ws.Range("E27").Value = "  +14.77%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "'2.527.55"
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("D30").Value = "'0.0₃0942"
$ws.Range("E30").Value = "  +4.86%  "
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("E32").Value = "  +5.63%  "
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("E34").Value = "  +3.34%  "
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("D36").Value = "'5.71"
$ws.Range("E36").Value = "  +4.22%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "'4.77"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("D39").Value = "'154.19"
$ws.Range("E39").Value = "  +5.12%  "
$ws.Range("E40").Value = "  +1.60%  "
$ws.Range("D41").Value = "'18.66"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("E42").Value = "  -4.75%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'2.36"
$ws.Range("E44").Value = "  +9.38%  "
$ws.Range("D45").Value = "'150.41"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("D47").Value = "'0.0538"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("D48").Value = "'20.35"
$ws.Range("E48").Value = "  +4.43%  "
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("E50").Value = "  +1.81%  "
